$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "34.139.98"
$ws.Cells.Item(2, 5).Value = "  +0.27%  "

$ws.Cells.Item(3, 4).Value = "1.792.75"
$ws.Cells.Item(3, 5).Value = "  +0.54%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.09%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "227.47"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.97%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.547"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -0.58%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -0.09%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "32.32"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -0.73%  "

$ws.Cells.Item(9, 5).Value = "  +3.24%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.0694"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.29%  "

$ws.Cells.Item(11, 5).Value = "  +0.67%  "

$ws.Cells.Item(12, 4).Value = "2.049.54"
$ws.Cells.Item(12, 5).Value = "  +0.44%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "11.57"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +5.12%  "

$ws.Cells.Item(14, 4).Value = "1.797.08"
$ws.Cells.Item(14, 5).Value = "  +0.62%  "

$ws.Cells.Item(15, 2).Value = "Polygon"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.624"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +0.74%  "

$ws.Cells.Item(16, 2).Value = "WrappedBTC"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Cells.Item(16, 4).Value = "34.125.71"
$ws.Cells.Item(16, 5).Value = "  +0.35%  "

$ws.Cells.Item(17, 5).Value = "  +1.32%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "68.02"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.31%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "245.81"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +0.53%  "

$ws.Cells.Item(20, 4).Value = "0.0₃0783"
$ws.Cells.Item(20, 5).Value = "  -0.14%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "10.94"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +2.07%  "

$ws.Cells.Item(22, 5).Value = "  +0.01%  "

$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +1.00%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "2.06"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.88%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "161.98"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.17%  "

$ws.Cells.Item(26, 5).Value = "  +2.09%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "16.32"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.08%  "

$ws.Cells.Item(28, 5).Value = "  +1.43%  "

$ws.Cells.Item(29, 5).Value = "  +0.00%  "

$ws.Cells.Item(30, 5).Value = "  +1.70%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "0.0521"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +1.66%  "

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "3.68"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +1.27%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "3.63"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +3.41%  "

$ws.Cells.Item(34, 5).Value = "  +1.25%  "

$ws.Cells.Item(35, 4).Value = "1.450.80"
$ws.Cells.Item(35, 5).Value = "  +4.34%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.648"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.40%  "

$ws.Cells.Item(37, 5).Value = "  +2.84%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "2.38"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +8.24%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "1.04"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -0.90%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "80.75"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  +3.45%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.929"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.06%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.35"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +0.23%  "

$ws.Cells.Item(43, 5).Value = "  +0.27%  "

$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "13.39"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +7.57%  "

$ws.Cells.Item(45, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(45, 4).Value = "0.0₆0140"
$ws.Cells.Item(45, 5).Value = "  -1.21%  "

$ws.Cells.Item(46, 2).Value = "FraxShare"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "6.07"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +4.20%  "

$ws.Cells.Item(47, 2).Value = "Kaspa"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.0509"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +2.34%  "

$ws.Cells.Item(48, 5).Value = "  -0.91%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "107.84"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.20%  "

$ws.Cells.Item(50, 4).Value = "1.951.21"
$ws.Cells.Item(50, 5).Value = "  +0.43%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -0.03%  "
